# Fruta / hortaliza, semanal
# A new weekly record was added to the "Apio" (celery) price series for
# "Vega Modelo de Temuco". In the canonical row order (sorted by date
# ascending) the new reading lands right after the current row 178, so it
# is inserted as row 179 and every following record (old rows 179-228)
# shifts down by one (new rows 180-229).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 179, pushing rows 179-228 down to 180-229.
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row 179 with the new weekly record.
$ws.Range("A179").Value = 10
$ws.Range("B179").Value = "Vega Modelo de Temuco"
$ws.Range("C179").Value = "La Araucanía"
$ws.Range("D179").Value = 44551
$ws.Range("E179").Value = 9
$ws.Range("F179").Value = 100112017
$ws.Range("G179").Value = "Apio"
$ws.Range("H179").Value = "Americana (o)"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 65
$ws.Range("K179").Value = 9000
$ws.Range("L179").Value = 9000
$ws.Range("M179").Value = 9000
$ws.Range("N179").Value = "`$/docena de matas"
$ws.Range("O179").Value = "Provincia del Elquí"
$ws.Range("P179").Value = 1500
$ws.Range("Q179").Value = 6
$ws.Range("R179").Value = "Hortaliza"

# Keep the date column formatted the same way as its neighbours.
$ws.Range("D179").NumberFormat = $ws.Range("D180").NumberFormat
